$wb = $excel.ActiveWorkbook

# Update the "Python_Trans_Data" sheet: replace formula-driven Capacity values
# with plain hardcoded values (150 instead of the formula that evaluated to 50)
$ws = $wb.Worksheets.Item("Python_Trans_Data")
$ws.Range("D2").Value = 150
$ws.Range("D3").Value = 150

# Move the active selection to D4 (as last edited cell) and make sure this
# sheet is the active/tab-selected one
$ws.Activate()
$ws.Range("D4").Select()
